$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) to stay text while we write new values, so Excel
# does not auto-convert plain-looking decimals (e.g. "69.35") into numbers.
# Column B/C/E values never look like bare numbers, so they are safe as-is.
$priceCol = $ws.Range("D2:D51")
$priceCol.NumberFormat = "@"

$ws.Range("D2").Value = '69.523.37'
$ws.Range("E2").Value = '  +0.80%  '

$ws.Range("D3").Value = '2.494.42'
$ws.Range("E3").Value = '  +0.12%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = '570.33'
$ws.Range("E5").Value = '  +0.30%  '

$ws.Range("D6").Value = '166.95'
$ws.Range("E6").Value = '  +1.61%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("E8").Value = '  -0.23%  '

$ws.Range("E9").Value = '  +0.91%  '

$ws.Range("E10").Value = '  -0.69%  '

$ws.Range("E11").Value = '  -0.58%  '

$ws.Range("E12").Value = '  +0.36%  '

$ws.Range("D13").Value = '2.950.19'
$ws.Range("E13").Value = '  -0.17%  '

$ws.Range("D14").Value = '69.369.52'
$ws.Range("E14").Value = '  +0.65%  '

$ws.Range("E15").Value = '  +1.15%  '

$ws.Range("D16").Value = '24.26'
$ws.Range("E16").Value = '  -1.20%  '

$ws.Range("D17").Value = '2.492.11'
$ws.Range("E17").Value = '  -0.58%  '

$ws.Range("D18").Value = '11.23'
$ws.Range("E18").Value = '  +0.51%  '

$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '7.47'
$ws.Range("E19").Value = '  -1.72%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '353.13'
$ws.Range("E20").Value = '  +2.21%  '

$ws.Range("D21").Value = '3.91'
$ws.Range("E21").Value = '  +1.08%  '

$ws.Range("E22").Value = '  -2.82%  '

$ws.Range("E23").Value = '  -0.12%  '

$ws.Range("D24").Value = '69.35'
$ws.Range("E24").Value = '  -0.79%  '

$ws.Range("D25").Value = '3.80'
$ws.Range("E25").Value = '  -1.84%  '

$ws.Range("D26").Value = '2.621.94'
$ws.Range("E26").Value = '  -1.03%  '

$ws.Range("D27").Value = '8.64'
$ws.Range("E27").Value = '  -1.48%  '

$ws.Range("D28").Value = '1.01'
$ws.Range("E28").Value = '  +0.53%  '

$ws.Range("D29").Value = '0.0₃0874'
$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("D30").Value = '7.59'
$ws.Range("E30").Value = '  -1.47%  '

$ws.Range("D31").Value = '443.82'
$ws.Range("E31").Value = '  -2.58%  '

$ws.Range("E32").Value = '  -1.47%  '

$ws.Range("B33").Value = 'POPCAT'
$ws.Range("C33").Value = 'https://coinranking.com/coin/sLBuDEsp6+popcat-popcat'
$ws.Range("D33").Value = '3.31'
$ws.Range("E33").Value = '  +119.75%  '

$ws.Range("E34").Value = '  +0.08%  '

$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.10%  '

$ws.Range("D36").Value = '154.46'
$ws.Range("E36").Value = '  -0.66%  '

$ws.Range("E37").Value = '  -1.38%  '

$ws.Range("E38").Value = '  +0.39%  '

$ws.Range("D39").Value = '18.15'
$ws.Range("E39").Value = '  -0.87%  '

$ws.Range("E40").Value = '  +0.04%  '

$ws.Range("E41").Value = '  -0.67%  '

$ws.Range("D42").Value = '4.63'
$ws.Range("E42").Value = '  -0.07%  '

$ws.Range("E43").Value = '  -0.31%  '

$ws.Range("E44").Value = '  +0.44%  '

$ws.Range("E45").Value = '  -3.52%  '

$ws.Range("D46").Value = '139.19'
$ws.Range("E46").Value = '  -1.19%  '

$ws.Range("D47").Value = '3.44'
$ws.Range("E47").Value = '  +0.15%  '

$ws.Range("D48").Value = '0.506'
$ws.Range("E48").Value = '  -1.30%  '

$ws.Range("D49").Value = '0.0724'
$ws.Range("E49").Value = '  -0.63%  '

$ws.Range("D50").Value = '0.572'
$ws.Range("E50").Value = '  -0.36%  '

$ws.Range("D51").Value = '0.0924'
$ws.Range("E51").Value = '  -0.05%  '

# Restore the default (General) style on column D now that the text values
# are locked in, so no stray number formatting lingers on the sheet.
$priceCol.Style = "Normal"
